$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.748.56"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "2.342.67"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "503.63"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "128.90"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "2.351.42"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "4.77"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "2.758.19"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "21.64"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "55.681.35"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "2.260.42"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").Value = "9.92"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("D20").Value = "310.25"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "3.99"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D26").Value = "0.371"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").Value = "0.146"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").Value = "7.08"
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").Value = "171.06"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "0.0₃0702"
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "5.76"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("D36").Value = "17.66"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "1.17"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").Value = "3.63"
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "36.03"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "126.52"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("D45").Value = "0.554"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "0.0890"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "237.54"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "0.0474"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "16.71"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "0.953"
$ws.Range("E51").Value = "  +0.02%  "
